$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '71.741.28'
$ws.Range('E2').Value = '  +4.73%  '
$ws.Range('D3').Value = '4.044.25'
$ws.Range('E3').Value = '  +4.78%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '533.50'
$ws.Range('E5').Value = '  +2.59%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '154.38'
$ws.Range('E6').Value = '  +9.53%  '
$ws.Range('E7').Value = '  +14.27%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.761'
$ws.Range('E9').Value = '  +7.08%  '
$ws.Range('E10').Value = '  +5.77%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0000333'
$ws.Range('E11').Value = '  +4.52%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '49.48'
$ws.Range('E12').Value = '  +19.22%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '11.00'
$ws.Range('E13').Value = '  +6.74%  '
$ws.Range('D14').Value = '4.688.59'
$ws.Range('E14').Value = '  +4.76%  '
$ws.Range('D15').Value = '4.040.81'
$ws.Range('E15').Value = '  +3.88%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '14.46'
$ws.Range('E16').Value = '  +2.26%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '21.10'
$ws.Range('E17').Value = '  -1.62%  '
$ws.Range('E18').Value = '  +2.17%  '
$ws.Range('E19').Value = '  -0.07%  '
$ws.Range('D20').Value = '71.786.07'
$ws.Range('E20').Value = '  +4.66%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '437.56'
$ws.Range('E21').Value = '  +5.21%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '3.72'
$ws.Range('E22').Value = '  +7.29%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '99.85'
$ws.Range('E23').Value = '  +15.20%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '14.86'
$ws.Range('E24').Value = '  +6.40%  '
$ws.Range('E25').Value = '  +6.71%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '11.43'
$ws.Range('E26').Value = '  +0.34%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.94'
$ws.Range('E27').Value = '  +4.11%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '37.36'
$ws.Range('E28').Value = '  +5.66%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '5.82'
$ws.Range('E29').Value = '  +2.85%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '3.54'
$ws.Range('E30').Value = '  +27.60%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '13.70'
$ws.Range('E31').Value = '  +3.99%  '
$ws.Range('E32').Value = '  +6.19%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '674.26'
$ws.Range('E33').Value = '  -0.44%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '6.76'
$ws.Range('E34').Value = '  +1.00%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '66.77'
$ws.Range('E35').Value = '  +2.33%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '43.21'
$ws.Range('E36').Value = '  +9.13%  '
$ws.Range('E37').Value = '  -2.54%  '
$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.158'
$ws.Range('E38').Value = '  +6.62%  '
$ws.Range('B39').Value = 'PEPE'
$ws.Range('C39').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D39').Value = '0.0₃0860'
$ws.Range('E39').Value = '  +3.68%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '3.43'
$ws.Range('E40').Value = '  -2.57%  '
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('E42').Value = '  +5.43%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.998'
$ws.Range('E43').Value = '  -0.22%  '
$ws.Range('E44').Value = '  +3.29%  '
$ws.Range('E45').Value = '  +9.19%  '
$ws.Range('E46').Value = '  -0.44%  '
$ws.Range('E47').Value = '  -0.31%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '9.54'
$ws.Range('E48').Value = '  +11.57%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '3.07'
$ws.Range('E49').Value = '  +3.39%  '
$ws.Range('B50').Value = 'LidoDAOToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '3.37'
$ws.Range('E50').Value = '  +3.26%  '
$ws.Range('B51').Value = 'FLOKI'
$ws.Range('C51').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.000274'
$ws.Range('E51').Value = '  +1.83%  '
